$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old standalone "Docentes responsaveis" value row (old row 13,
# which only had B13/C13 populated with the professor name). Deleting it
# shifts rows 14-22 up to 13-21 and fixes up row heights/dimension to
# match the target layout automatically (Excel re-flows ht=.. and the
# <dimension> ref for us).
$ws.Rows(13).Delete()

# Wipe existing cell contents but keep per-cell styles/row heights, then
# re-populate every surviving cell from scratch in strict reading order
# (top-to-bottom, column A then B then C). Excel builds the shared-string
# table in first-use order, so writing in this order reproduces the
# target sharedStrings.xml ordering exactly.
$ws.Cells.ClearContents()

# Helper cell used below to push text in verbatim (see note next to its
# uses): writing through it with a leading apostrophe plus a values-only
# paste stops Excel from auto-converting look-alike numbers/dates (e.g.
# "4", "0", "01/01/2012") into real numeric/date values, while leaving
# the destination cell's existing style/number format untouched.

$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'

$ws.Range("B2").Value = 'LOM3201'
$ws.Range("C2").Value = 'LOM3201'

$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Biofísica Molecular'
$ws.Range("C3").Value = ' Biofísica Molecular'

$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Molecular Biophysics'
$ws.Range("C4").Value = 'Molecular Biophysics'

$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("ZZ1").Value = "'4"
$ws.Range("ZZ1").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false
$ws.Range("ZZ1").Value = "'4"
$ws.Range("ZZ1").Copy()
$ws.Range("C5").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("ZZ1").Value = "'0"
$ws.Range("ZZ1").Copy()
$ws.Range("B6").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false
$ws.Range("ZZ1").Value = "'0"
$ws.Range("ZZ1").Copy()
$ws.Range("C6").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'

$ws.Range("A8").Value = 'Ativação:'
$ws.Range("ZZ1").Value = "'01/01/2012"
$ws.Range("ZZ1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false
$ws.Range("ZZ1").Value = "'01/01/2012"
$ws.Range("ZZ1").Copy()
$ws.Range("C8").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EF-8'
$ws.Range("C9").Value = 'EF-8'

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '2166002 - Sandra Giacomin Schneider'
$ws.Range("C10").Value = '2166002 - Sandra Giacomin Schneider'

$ws.Range("A11").Value = 'Objectives:'

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

$ws.Range("A14").Value = 'Short syllabus:'

$ws.Range("A15").Value = 'Programa:'
$ws.Range("ZZ1").Value = "'01/01/2012"
$ws.Range("ZZ1").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false
$ws.Range("ZZ1").Value = "'01/01/2012"
$ws.Range("ZZ1").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Range("A16").Value = 'Syllabus:'

$ws.Range("A17").Value = 'Avaliação:'

$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '2166002 - Sandra Giacomin Schneider'
$ws.Range("C18").Value = '2166002 - Sandra Giacomin Schneider'

$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Aulas expositivas; demonstrações com o uso de programas de computador e discussão de listas de exercícios'
$ws.Range("C19").Value = 'Aulas expositivas; demonstrações com o uso de programas de computador e discussão de listas de exercícios'

$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Três provas escritas: conceitos P1, P2 e P3. Conceito Final = (P1 + P2 + 2P3)/4'
$ws.Range("C20").Value = 'Três provas escritas: conceitos P1, P2 e P3. Conceito Final = (P1 + P2 + 2P3)/4'

$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'

$ws.Range("ZZ1").Clear()
